$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.453.29'
$ws.Range('E2').Value = '  -2.43%  '
$ws.Range('D3').Value = '3.698.97'
$ws.Range('E3').Value = '  -3.13%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '693.01'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '163.40'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.17%  '
$ws.Range('D7').Value = '3.697.20'
$ws.Range('E7').Value = '  -3.15%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -4.27%  '
$ws.Range('E10').Value = '  -7.90%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.39'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.71%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.446'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000240'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -5.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '33.51'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -6.89%  '
$ws.Range('D15').Value = '4.322.61'
$ws.Range('E15').Value = '  -3.16%  '
$ws.Range('D16').Value = '3.703.39'
$ws.Range('E16').Value = '  -3.05%  '
$ws.Range('D17').Value = '69.507.57'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '16.33'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -6.98%  '
$ws.Range('E20').Value = '  -7.47%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '482.23'
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.99'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -6.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.667'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -7.60%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '80.15'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -4.76%  '
$ws.Range('D25').Value = '3.843.60'
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('E26').Value = '  -9.19%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.44'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.16%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.55'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -8.41%  '
$ws.Range('E30').Value = '  -9.63%  '
$ws.Range('E31').Value = '  -9.64%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.88'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -7.16%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.08'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -7.24%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '27.13'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.75%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.167'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.97%  '
$ws.Range('D37').Value = '3.665.95'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.52'
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.39'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +6.38%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.34'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.89%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0936'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -7.40%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.953'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -6.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '163.91'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '48.08'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.80%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '30.12'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.83'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -14.15%  '
$ws.Range('B49').Value = 'SuiNetwork'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.16'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('E51').Value = '  -7.57%  '
